$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Physiology")
$ws.Rows.Item(10).Insert()
$ws.Range("A10:I10").RowHeight = 84.75
$ws.Range("B10").HorizontalAlignment = -4108
Write-Output "done"
